$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix hierarchical metrics wrongly computed
# Column C (btop) corrected values for rows 2-4
$ws.Range("C2").Value = 0.8471901075701217
$ws.Range("C3").Value = 0.8471901075701217
$ws.Range("C4").Value = 0.8471901075701218

# Column H (hitac_filter_qiime) corrected values for rows 2-4 (now equal to column I)
$ws.Range("H2").Value = 0.93711467324291
$ws.Range("H3").Value = 0.8410462776659959
$ws.Range("H4").Value = 0.8864853401198238
